$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (Tipo) so MAE becomes column D and Tipo shifts to E
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("D1").Value = "MAE"
$ws.Range("E1").Value = "Tipo"

# Data row - updated values
$ws.Range("B2").Value = 0.2844791062672977
$ws.Range("C2").Value = 0.9945272545749643
$ws.Range("D2").Value = 0.4311542036564542
$ws.Range("E2").Value = "single"
